$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" headers (A1:J1) to "_FV2310" and the "_new" headers
# (L1:U1) to "_FV2404" -- column K1 ("diff") is left untouched.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $cell.Value().Replace("_old", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $cell.Value().Replace("_new", "_FV2404")
}

# Turn the data range into a proper Excel Table (ListObject).
$rng = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
